# Removed date conversion in ClaimProcessMarine
# Fix the "CauseofLoss" / "eBaoClass" mapping labels on the Mapping sheet:
#   C17: "Casue Of Loss" (typo)  -> "Cause Of Loss"
#   C19: "eBao"                 -> "eBao Class"
# and move the active selection from C20 to C19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("C17").Value = "Cause Of Loss"
$ws.Range("C19").Value = "eBao Class"

$ws.Activate()
$ws.Range("C19").Select()
